# Auto-generated edit script: update market-price-derived columns (H-N)
# per sheet/row, matching the scheduled market-data refresh described in the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 2507.0715
$ws.Cells.Item(12, 9).Value = 1566.5555
$ws.Cells.Item(12, 11).Value = 1566.5555
$ws.Cells.Item(12, 13).Value = -1396.5555

$ws.Cells.Item(17, 8).Value = 2489.1052
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 2489.1052
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 7467.3156
$ws.Cells.Item(17, 14).Value = -7803.3156
$ws.Cells.Item(17, 13).ClearContents()

$ws.Cells.Item(57, 8).Value = 168386.72
$ws.Cells.Item(57, 10).Value = 168386.72
$ws.Cells.Item(57, 12).Value = 505160.16
$ws.Cells.Item(57, 14).Value = -506158.16

$ws.Cells.Item(86, 8).Value = 2294.9614
$ws.Cells.Item(86, 9).Value = 2065.8
$ws.Cells.Item(86, 10).Value = 2607.4546
$ws.Cells.Item(86, 11).Value = 2065.8
$ws.Cells.Item(86, 12).Value = 2607.4546
$ws.Cells.Item(86, 13).Value = -942.8000000000002
$ws.Cells.Item(86, 14).Value = -4853.4546

$ws.Cells.Item(89, 8).Value = 2294.9614
$ws.Cells.Item(89, 9).Value = 2065.8
$ws.Cells.Item(89, 10).Value = 2607.4546
$ws.Cells.Item(89, 11).Value = 10329
$ws.Cells.Item(89, 12).Value = 13037.273
$ws.Cells.Item(89, 13).Value = -4713
$ws.Cells.Item(89, 14).Value = -24269.273

$ws.Cells.Item(112, 8).Value = 1115437.9
$ws.Cells.Item(112, 9).Value = 1443
$ws.Cells.Item(112, 10).Value = 1254687.2
$ws.Cells.Item(112, 11).Value = 4329
$ws.Cells.Item(112, 12).Value = 3764061.6
$ws.Cells.Item(112, 13).Value = -3221
$ws.Cells.Item(112, 14).Value = -3766277.6

$ws.Cells.Item(116, 8).Value = 2700
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(125, 8).Value = 6108.778
$ws.Cells.Item(125, 10).Value = 6108.778
$ws.Cells.Item(125, 12).Value = 54979.002
$ws.Cells.Item(125, 14).Value = -59899.002

$ws.Cells.Item(137, 8).Value = 2566.3901
$ws.Cells.Item(137, 9).Value = 2392.5715
$ws.Cells.Item(137, 10).Value = 3580.3333
$ws.Cells.Item(137, 11).Value = 7177.7145
$ws.Cells.Item(137, 12).Value = 10740.9999
$ws.Cells.Item(137, 13).Value = -4627.7145
$ws.Cells.Item(137, 14).Value = -15840.9999

$ws.Cells.Item(138, 8).Value = 381474.94
$ws.Cells.Item(138, 10).Value = 718010.3
$ws.Cells.Item(138, 12).Value = 2154030.9
$ws.Cells.Item(138, 14).Value = -2164310.9

$ws = $wb.Sheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6898.6797
$ws.Cells.Item(32, 9).Value = 6931.04
$ws.Cells.Item(32, 10).Value = 6089.6665
$ws.Cells.Item(32, 11).Value = 6931.04
$ws.Cells.Item(32, 12).Value = 6089.6665
$ws.Cells.Item(32, 13).Value = -6644.04
$ws.Cells.Item(32, 14).Value = -6663.6665

$ws.Cells.Item(45, 8).Value = 3307.4
$ws.Cells.Item(45, 9).Value = 2561.2
$ws.Cells.Item(45, 10).Value = 4799.8
$ws.Cells.Item(45, 11).Value = 2561.2
$ws.Cells.Item(45, 12).Value = 4799.8
$ws.Cells.Item(45, 13).Value = -2184.2
$ws.Cells.Item(45, 14).Value = -5553.8

$ws.Cells.Item(61, 8).Value = 3693.0193
$ws.Cells.Item(61, 9).Value = 3982.818
$ws.Cells.Item(61, 10).Value = 2099.125
$ws.Cells.Item(61, 11).Value = 3982.818
$ws.Cells.Item(61, 12).Value = 2099.125
$ws.Cells.Item(61, 13).Value = -3770.818
$ws.Cells.Item(61, 14).Value = -2523.125

$ws.Cells.Item(122, 8).Value = 1250.1111
$ws.Cells.Item(122, 9).Value = 1144.3462
$ws.Cells.Item(122, 11).Value = 3433.0386
$ws.Cells.Item(122, 13).Value = -983.0385999999999

$ws.Cells.Item(132, 8).Value = 3835.74
$ws.Cells.Item(132, 9).Value = 2182.244
$ws.Cells.Item(132, 11).Value = 6546.732
$ws.Cells.Item(132, 13).Value = -4016.732

$ws.Cells.Item(136, 8).Value = 3693.0193
$ws.Cells.Item(136, 9).Value = 3982.818
$ws.Cells.Item(136, 10).Value = 2099.125
$ws.Cells.Item(136, 11).Value = 11948.454
$ws.Cells.Item(136, 12).Value = 6297.375
$ws.Cells.Item(136, 13).Value = -9398.454000000002
$ws.Cells.Item(136, 14).Value = -11397.375

$ws = $wb.Sheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 735.875
$ws.Cells.Item(80, 9).Value = 588.8
$ws.Cells.Item(80, 10).Value = 802.7273
$ws.Cells.Item(80, 11).Value = 588.8
$ws.Cells.Item(80, 12).Value = 802.7273
$ws.Cells.Item(80, 13).Value = 409.2
$ws.Cells.Item(80, 14).Value = -2798.7273

$ws.Cells.Item(83, 8).Value = 735.875
$ws.Cells.Item(83, 9).Value = 588.8
$ws.Cells.Item(83, 10).Value = 802.7273
$ws.Cells.Item(83, 11).Value = 2944
$ws.Cells.Item(83, 12).Value = 4013.6365
$ws.Cells.Item(83, 13).Value = 2048
$ws.Cells.Item(83, 14).Value = -13997.6365

$ws.Cells.Item(94, 8).Value = 1842.85
$ws.Cells.Item(94, 9).Value = 1621.5883
$ws.Cells.Item(94, 11).Value = 1621.5883
$ws.Cells.Item(94, 13).Value = -1170.5883

$ws.Cells.Item(105, 8).Value = 2844.6453
$ws.Cells.Item(105, 9).Value = 2218
$ws.Cells.Item(105, 10).Value = 5455.6665
$ws.Cells.Item(105, 11).Value = 2218
$ws.Cells.Item(105, 12).Value = 5455.6665
$ws.Cells.Item(105, 13).Value = -471
$ws.Cells.Item(105, 14).Value = -8949.666499999999

$ws = $wb.Sheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1782.6538
$ws.Cells.Item(31, 9).Value = 1382.289
$ws.Cells.Item(31, 11).Value = 1382.289
$ws.Cells.Item(31, 13).Value = -1087.289

$ws.Cells.Item(34, 8).Value = 1782.6538
$ws.Cells.Item(34, 9).Value = 1382.289
$ws.Cells.Item(34, 11).Value = 1382.289
$ws.Cells.Item(34, 13).Value = -1180.289

$ws.Cells.Item(86, 8).Value = 7174.5
$ws.Cells.Item(86, 9).Value = 7066
$ws.Cells.Item(86, 11).Value = 7066
$ws.Cells.Item(86, 13).Value = -5943

$ws.Cells.Item(89, 8).Value = 7174.5
$ws.Cells.Item(89, 9).Value = 7066
$ws.Cells.Item(89, 11).Value = 35330
$ws.Cells.Item(89, 13).Value = -29714

$ws.Cells.Item(94, 8).Value = 1217.0476
$ws.Cells.Item(94, 10).Value = 1392
$ws.Cells.Item(94, 12).Value = 1392
$ws.Cells.Item(94, 14).Value = -2294

$ws.Cells.Item(132, 8).Value = 2300.7334
$ws.Cells.Item(132, 9).Value = 2126.1765
$ws.Cells.Item(132, 11).Value = 6378.529500000001
$ws.Cells.Item(132, 13).Value = -3848.529500000001

$ws = $wb.Sheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 50.625
$ws.Cells.Item(2, 9).Value = 60.5
$ws.Cells.Item(2, 10).Value = 30.875
$ws.Cells.Item(2, 11).Value = 363
$ws.Cells.Item(2, 12).Value = 185.25
$ws.Cells.Item(2, 13).Value = -250
$ws.Cells.Item(2, 14).Value = -411.25

$ws.Cells.Item(40, 8).Value = 363.08334
$ws.Cells.Item(40, 9).Value = 51
$ws.Cells.Item(40, 10).Value = 800
$ws.Cells.Item(40, 11).Value = 204
$ws.Cells.Item(40, 12).Value = 3200
$ws.Cells.Item(40, 13).Value = -135
$ws.Cells.Item(40, 14).Value = -3338

$ws.Cells.Item(87, 8).Value = 1679839.4
$ws.Cells.Item(87, 9).Value = 1679839.4
$ws.Cells.Item(87, 11).Value = 5039518.199999999
$ws.Cells.Item(87, 13).Value = -5038270.199999999

$ws.Cells.Item(90, 8).Value = 1679839.4
$ws.Cells.Item(90, 9).Value = 1679839.4
$ws.Cells.Item(90, 11).Value = 15118554.6
$ws.Cells.Item(90, 13).Value = -15112314.6

$ws.Cells.Item(107, 8).Value = 2176.8064
$ws.Cells.Item(107, 9).Value = 188.6
$ws.Cells.Item(107, 10).Value = 3123.5715
$ws.Cells.Item(107, 11).Value = 565.8
$ws.Cells.Item(107, 12).Value = 9370.7145
$ws.Cells.Item(107, 13).Value = 1354.2
$ws.Cells.Item(107, 14).Value = -13210.7145

$ws = $wb.Sheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 58195.59
$ws.Cells.Item(43, 10).Value = 77908.09
$ws.Cells.Item(43, 12).Value = 77908.09
$ws.Cells.Item(43, 14).Value = -78210.09

$ws.Cells.Item(126, 8).Value = 2740.8235
$ws.Cells.Item(126, 9).Value = 2537.75
$ws.Cells.Item(126, 10).Value = 3228.2
$ws.Cells.Item(126, 11).Value = 7613.25
$ws.Cells.Item(126, 12).Value = 9684.599999999999
$ws.Cells.Item(126, 13).Value = -5143.25
$ws.Cells.Item(126, 14).Value = -14624.6

$ws.Cells.Item(132, 8).Value = 4059.55
$ws.Cells.Item(132, 9).Value = 4121.8613
$ws.Cells.Item(132, 11).Value = 12365.5839
$ws.Cells.Item(132, 13).Value = -9835.583899999998

$ws = $wb.Sheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 9963.4
$ws.Cells.Item(46, 9).Value = 3327.1428
$ws.Cells.Item(46, 10).Value = 15770.125
$ws.Cells.Item(46, 11).Value = 3327.1428
$ws.Cells.Item(46, 12).Value = 15770.125
$ws.Cells.Item(46, 13).Value = -3139.1428
$ws.Cells.Item(46, 14).Value = -16146.125

$ws.Cells.Item(122, 8).Value = 5345.8096
$ws.Cells.Item(122, 9).Value = 4082.077
$ws.Cells.Item(122, 11).Value = 12246.231
$ws.Cells.Item(122, 13).Value = -9796.231

$ws.Cells.Item(132, 8).Value = 2225.3555
$ws.Cells.Item(132, 9).Value = 1332.5938
$ws.Cells.Item(132, 10).Value = 4422.923
$ws.Cells.Item(132, 11).Value = 3997.7814
$ws.Cells.Item(132, 12).Value = 13268.769
$ws.Cells.Item(132, 13).Value = -1467.7814
$ws.Cells.Item(132, 14).Value = -18328.769

$ws.Cells.Item(136, 8).Value = 4085.465
$ws.Cells.Item(136, 9).Value = 3623.28
$ws.Cells.Item(136, 10).Value = 4727.3887
$ws.Cells.Item(136, 11).Value = 10869.84
$ws.Cells.Item(136, 12).Value = 14182.1661
$ws.Cells.Item(136, 13).Value = -8319.84
$ws.Cells.Item(136, 14).Value = -19282.1661

$ws = $wb.Sheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 12600717
$ws.Cells.Item(62, 9).Value = 133457.67
$ws.Cells.Item(62, 11).Value = 133457.67
$ws.Cells.Item(62, 13).Value = -132833.67

$ws.Cells.Item(65, 8).Value = 12600717
$ws.Cells.Item(65, 9).Value = 133457.67
$ws.Cells.Item(65, 11).Value = 667288.3500000001
$ws.Cells.Item(65, 13).Value = -664168.3500000001

$ws.Cells.Item(76, 8).Value = 24000
$ws.Cells.Item(76, 10).Value = 24000
$ws.Cells.Item(76, 12).Value = 24000
$ws.Cells.Item(76, 14).Value = -24630

$ws.Cells.Item(79, 8).Value = 24000
$ws.Cells.Item(79, 10).Value = 24000
$ws.Cells.Item(79, 12).Value = 24000
$ws.Cells.Item(79, 14).Value = -26184

$ws.Cells.Item(81, 8).Value = 39002.1
$ws.Cells.Item(81, 9).Value = 74373.07000000001
$ws.Cells.Item(81, 11).Value = 148746.14
$ws.Cells.Item(81, 13).Value = -147685.14

$ws.Cells.Item(84, 8).Value = 39002.1
$ws.Cells.Item(84, 9).Value = 74373.07000000001
$ws.Cells.Item(84, 11).Value = 743730.7000000001
$ws.Cells.Item(84, 13).Value = -738426.7000000001

$ws.Cells.Item(100, 8).Value = 1582.0435
$ws.Cells.Item(100, 9).Value = 1534.4706
$ws.Cells.Item(100, 11).Value = 3068.9412
$ws.Cells.Item(100, 13).Value = -2527.9412

$ws.Cells.Item(122, 8).Value = 2453.3684
$ws.Cells.Item(122, 9).Value = 2104.1292
$ws.Cells.Item(122, 10).Value = 4000
$ws.Cells.Item(122, 11).Value = 6312.3876
$ws.Cells.Item(122, 12).Value = 12000
$ws.Cells.Item(122, 13).Value = -3862.3876
$ws.Cells.Item(122, 14).Value = -16900

$ws.Cells.Item(136, 8).Value = 13578.258
$ws.Cells.Item(136, 9).Value = 18437.475
$ws.Cells.Item(136, 11).Value = 55312.425
$ws.Cells.Item(136, 13).Value = -52762.425
